$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $r = $Sheet.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.ClearFormats()
}

Set-TextValue $ws "D2" "27.527.54"
Set-TextValue $ws "E2" "  +5.31%  "

Set-TextValue $ws "D3" "1.724.51"
Set-TextValue $ws "E3" "  +4.12%  "

Set-TextValue $ws "E4" "  +0.08%  "

Set-TextValue $ws "D5" "225.95"
Set-TextValue $ws "E5" "  +3.36%  "

Set-TextValue $ws "D6" "0.5366"
Set-TextValue $ws "E6" "  +2.47%  "

Set-TextValue $ws "D8" "0.2676"
Set-TextValue $ws "E8" "  +0.60%  "

Set-TextValue $ws "D9" "0.06599"
Set-TextValue $ws "E9" "  +3.81%  "

Set-TextValue $ws "D10" "21.74"
Set-TextValue $ws "E10" "  +5.63%  "

Set-TextValue $ws "E11" "  +0.76%  "

Set-TextValue $ws "D12" "4.624"
Set-TextValue $ws "E12" "  +0.27%  "

Set-TextValue $ws "D13" "1.719.81"
Set-TextValue $ws "E13" "  +6.09%  "

Set-TextValue $ws "D14" "1.962.24"
Set-TextValue $ws "E14" "  +4.15%  "

Set-TextValue $ws "D15" "0.5867"
Set-TextValue $ws "E15" "  +4.27%  "

Set-TextValue $ws "D16" "0.0₅8306"
Set-TextValue $ws "E16" "  +1.19%  "

Set-TextValue $ws "D17" "68.04"
Set-TextValue $ws "E17" "  +3.85%  "

Set-TextValue $ws "D18" "27.545.70"
Set-TextValue $ws "E18" "  +5.35%  "

Set-TextValue $ws "D19" "222.94"
Set-TextValue $ws "E19" "  +15.64%  "

Set-TextValue $ws "D20" "1.003"
Set-TextValue $ws "E20" "  +0.01%  "

Set-TextValue $ws "D22" "10.68"
Set-TextValue $ws "E22" "  +1.40%  "

Set-TextValue $ws "D23" "6.100"
Set-TextValue $ws "E23" "  +2.41%  "

Set-TextValue $ws "E24" "  +0.02%  "

Set-TextValue $ws "D25" "147.95"
Set-TextValue $ws "E25" "  +1.71%  "

Set-TextValue $ws "D26" "1.693"
Set-TextValue $ws "E26" "  +12.29%  "

Set-TextValue $ws "D27" "0.1231"
Set-TextValue $ws "E27" "  +2.84%  "

Set-TextValue $ws "D28" "7.392"
Set-TextValue $ws "E28" "  +1.71%  "

Set-TextValue $ws "D29" "16.68"
Set-TextValue $ws "E29" "  +4.43%  "

Set-TextValue $ws "D30" "0.05537"
Set-TextValue $ws "E30" "  +1.31%  "

Set-TextValue $ws "E31" "  +2.51%  "

Set-TextValue $ws "D32" "3.547"
Set-TextValue $ws "E32" "  +2.33%  "

Set-TextValue $ws "D33" "3.463"
Set-TextValue $ws "E33" "  +2.63%  "

Set-TextValue $ws "D34" "1.662"
Set-TextValue $ws "E34" "  +6.03%  "

Set-TextValue $ws "D35" "0.9597"
Set-TextValue $ws "E35" "  +0.61%  "

Set-TextValue $ws "B36" "MXToken"
Set-TextValue $ws "C36" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws "D36" "2.821"
Set-TextValue $ws "E36" "  +1.52%  "

Set-TextValue $ws "B37" "HuobiToken"
Set-TextValue $ws "C37" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws "D37" "2.445"
Set-TextValue $ws "E37" "  +1.77%  "

Set-TextValue $ws "D38" "0.5941"

Set-TextValue $ws "D39" "0.01649"
Set-TextValue $ws "E39" "  +3.88%  "

Set-TextValue $ws "D40" "5.908"
Set-TextValue $ws "E40" "  +0.49%  "

Set-TextValue $ws "D41" "1.060.28"
Set-TextValue $ws "E41" "  +3.31%  "

Set-TextValue $ws "D42" "0.8566"
Set-TextValue $ws "E42" "  +2.88%  "

Set-TextValue $ws "D44" "101.69"
Set-TextValue $ws "E44" "  +0.35%  "

Set-TextValue $ws "D45" "1.868.12"
Set-TextValue $ws "E45" "  +4.06%  "

Set-TextValue $ws "D46" "0.0₈115"
Set-TextValue $ws "E46" "  +10.46%  "

Set-TextValue $ws "D47" "59.04"
Set-TextValue $ws "E47" "  +2.15%  "

Set-TextValue $ws "D48" "8.212"
Set-TextValue $ws "E48" "  +2.60%  "

Set-TextValue $ws "D50" "1.005"
Set-TextValue $ws "E50" "  +0.44%  "

Set-TextValue $ws "D51" "0.05278"
Set-TextValue $ws "E51" "  +1.54%  "

